$wb = $excel.ActiveWorkbook

# Sheet "展览" : F2 255 -> 256, F3 370 -> 371
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 256
$ws1.Range("F3").Value = 371

# Sheet "全部类型" : F2 255 -> 256, F3 370 -> 371
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 256
$ws4.Range("F3").Value = 371
